# Recalibrate the model with obs till 2018-06-30 and update figures.
#
# The underlying parameter table is unchanged in value/meaning; the two
# "ScalingFactor (...)" labels in column I (rows 2 and 3) are renamed to
# "ScaleFactor (...)" to match the updated naming used elsewhere in the
# recalibrated model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "ScaleFactor (Subsurf)"
$ws.Range("I3").Value = "ScaleFactor (Surf)"

# Scroll the view so column C becomes the left-most visible column,
# without disturbing the current selection (stays on I12).
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
